$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The uploaded workbook replaced every occurrence of the text "n" (an
# unused/meaningless placeholder string) in column E with the number 6,
# which also drops "n" from the shared-strings table once nothing else
# references it.
$rows = @(4, 8, 10, 14, 20, 22, 25, 31, 35, 39, 41, 45)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = 6
}

# The saved view's active cell/selection moved from C41 to E39.
$null = $ws.Range("E39").Select()
